# Generate Report for Handoff
# - Priority for rows 4-7 (the 217254bf-... file entries) changes from "low" to "ht"
#   on both the zh-cn and de-de status sheets.
# - The "Latest Handoff Datetime" for those same rows is refreshed:
#     zh-cn sheet: 2016-10-27 10:54:03 -> 2016-10-27 10:54:57
#     de-de sheet: 2016-10-27 10:54:17 -> 2016-10-27 10:55:12

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

foreach ($row in 4..7) {
    $zh.Cells.Item($row, 5).Value = "ht"
    $de.Cells.Item($row, 5).Value = "ht"

    $zh.Cells.Item($row, 8).Value = "2016-10-27 10:54:57"
    $de.Cells.Item($row, 8).Value = "2016-10-27 10:55:12"
}
